# Apply the Golden State Warriors 2025-26 stats update:
#  1. Append 2025-11-11 (OKC) and 2025-11-12 (SAS) game rows to the
#     Points / Assists / Rebounds / 3PM sheets.
#  2. Recompute (re-write) the Avg Points / Avg Assists / Avg Rebounds /
#     Avg 3PM sheets with the new per-player averages, re-sorted
#     descending by average.
#  3. Add a new "Team Points" sheet with team/opponent/total scoring
#     for every game played so far, including the two new games.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Append the two new game rows to each boxscore sheet.
# ---------------------------------------------------------------------

# Writes a "yyyy-mm-dd" looking string into a cell while keeping it a
# plain text value (Excel otherwise auto-converts such strings into
# date values). Restoring the Normal style afterwards keeps the cell
# free of an explicit number-format override, just like its neighbors.
function Set-TextValue {
    param($cell, $text)

    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Add-GameRows {
    param($sheetName, $row13, $row14)

    $ws = $wb.Worksheets.Item($sheetName)

    Set-TextValue $ws.Cells.Item(13, 1) $row13[0]
    $ws.Cells.Item(13, 2).Value = $row13[1]
    for ($i = 0; $i -lt 14; $i++) {
        $ws.Cells.Item(13, 3 + $i).Value = $row13[2 + $i]
    }

    Set-TextValue $ws.Cells.Item(14, 1) $row14[0]
    $ws.Cells.Item(14, 2).Value = $row14[1]
    for ($i = 0; $i -lt 14; $i++) {
        $ws.Cells.Item(14, 3 + $i).Value = $row14[2 + $i]
    }
}

# Columns C..P on every boxscore sheet are, in order:
# Gary Payton II, Jonathan Kuminga, Brandin Podziemski, Will Richard,
# Moses Moody, Buddy Hield, Jimmy Butler III, Gui Santos, Al Horford,
# Quinten Post, Draymond Green, Stephen Curry, Trayce Jackson-Davis,
# Pat Spencer

Add-GameRows "Points" `
    @("2025-11-11", "OKC", 2, 13, 10, 9, 10, 4, 12, 4, 0, 6, 3, 11, 6, 12) `
    @("2025-11-12", "SAS", 4, 0, 6, 3, 19, 4, 28, 0, 9, 0, 6, 46, 0, 0)

Add-GameRows "Assists" `
    @("2025-11-11", "OKC", 3, 4, 1, 1, 0, 4, 2, 0, 0, 2, 4, 0, 1, 1) `
    @("2025-11-12", "SAS", 2, 0, 4, 1, 3, 0, 8, 0, 3, 0, 4, 5, 0, 0)

Add-GameRows "Rebounds" `
    @("2025-11-11", "OKC", 6, 3, 4, 0, 0, 0, 3, 3, 0, 3, 2, 1, 6, 1) `
    @("2025-11-12", "SAS", 6, 4, 6, 1, 3, 1, 6, 0, 3, 0, 5, 5, 0, 0)

Add-GameRows "3PM" `
    @("2025-11-11", "OKC", 0, 0, 2, 3, 2, 0, 0, 0, 0, 2, 1, 1, 0, 2) `
    @("2025-11-12", "SAS", 0, 0, 2, 1, 5, 0, 5, 0, 3, 0, 0, 5, 0, 0)

# ---------------------------------------------------------------------
# 2. Rewrite the four "Avg ..." sheets with updated, re-sorted values.
# ---------------------------------------------------------------------

function Set-AvgSheet {
    param($sheetName, $rows)

    $ws = $wb.Worksheets.Item($sheetName)

    # Clear out the old player rows (keep the header in row 1).
    $ws.Range("A2:B15").ClearContents()

    $r = 2
    foreach ($pair in $rows) {
        $ws.Cells.Item($r, 1).Value = $pair[0]
        $ws.Cells.Item($r, 2).Value = $pair[1]
        $r += 1
    }
}

Set-AvgSheet "Avg Points" @(
    @("Stephen Curry", 27.1),
    @("Jimmy Butler III", 19.08333333333333),
    @("Jonathan Kuminga", 13.76923076923077),
    @("Moses Moody", 12.27272727272727),
    @("Brandin Podziemski", 11.53846153846154),
    @("Will Richard", 9.363636363636363),
    @("Draymond Green", 7.75),
    @("Buddy Hield", 7.153846153846154),
    @("Al Horford", 6),
    @("Quinten Post", 5.846153846153846),
    @("Pat Spencer", 5.125),
    @("Trayce Jackson-Davis", 4.5),
    @("Gary Payton II", 2.25),
    @("Gui Santos", 1.909090909090909)
)

Set-AvgSheet "Avg Assists" @(
    @("Draymond Green", 5.5),
    @("Jimmy Butler III", 4.666666666666667),
    @("Stephen Curry", 3.9),
    @("Brandin Podziemski", 3.153846153846154),
    @("Jonathan Kuminga", 2.846153846153846),
    @("Pat Spencer", 1.75),
    @("Moses Moody", 1.636363636363636),
    @("Buddy Hield", 1.615384615384615),
    @("Al Horford", 1.5),
    @("Will Richard", 1.454545454545455),
    @("Trayce Jackson-Davis", 1.25),
    @("Gary Payton II", 1.083333333333333),
    @("Quinten Post", 1),
    @("Gui Santos", 0.1818181818181818)
)

Set-AvgSheet "Avg Rebounds" @(
    @("Jonathan Kuminga", 6.615384615384615),
    @("Jimmy Butler III", 5.333333333333333),
    @("Draymond Green", 5.25),
    @("Brandin Podziemski", 4.769230769230769),
    @("Al Horford", 4.25),
    @("Quinten Post", 3.769230769230769),
    @("Stephen Curry", 3.5),
    @("Moses Moody", 2.636363636363636),
    @("Gary Payton II", 2.583333333333333),
    @("Trayce Jackson-Davis", 2.5),
    @("Will Richard", 2.181818181818182),
    @("Pat Spencer", 1.625),
    @("Gui Santos", 1.454545454545455),
    @("Buddy Hield", 1.384615384615385)
)

Set-AvgSheet "Avg 3PM" @(
    @("Stephen Curry", 4.1),
    @("Moses Moody", 2.818181818181818),
    @("Brandin Podziemski", 1.692307692307692),
    @("Draymond Green", 1.583333333333333),
    @("Al Horford", 1.5),
    @("Will Richard", 1.363636363636364),
    @("Quinten Post", 1.230769230769231),
    @("Jimmy Butler III", 1.166666666666667),
    @("Buddy Hield", 1.153846153846154),
    @("Jonathan Kuminga", 0.9230769230769231),
    @("Pat Spencer", 0.5),
    @("Gui Santos", 0.3636363636363636),
    @("Gary Payton II", 0.08333333333333333),
    @("Trayce Jackson-Davis", 0)
)

# ---------------------------------------------------------------------
# 3. Add the new "Team Points" sheet at the end of the workbook.
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$teamSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$teamSheet.Name = "Team Points"

$headers = @("Game Time (PST)", "Opponent", "Team Points", "Opponent Points", "Game Total Points")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $teamSheet.Cells.Item(1, $c + 1).Value = $headers[$c]
}
$headerRange = $teamSheet.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$teamGames = @(
    @("2025-10-21", "LAL", 119, 109, 228),
    @("2025-10-23", "DEN", 137, 131, 268),
    @("2025-10-24", "POR", 119, 139, 258),
    @("2025-10-27", "MEM", 131, 118, 249),
    @("2025-10-28", "LAC", 98, 79, 177),
    @("2025-10-30", "MIL", 110, 120, 230),
    @("2025-11-01", "IND", 109, 114, 223),
    @("2025-11-04", "PHX", 118, 107, 225),
    @("2025-11-05", "SAC", 116, 121, 237),
    @("2025-11-07", "DEN", 104, 129, 233),
    @("2025-11-09", "IND", 114, 83, 197),
    @("2025-11-11", "OKC", 102, 126, 228),
    @("2025-11-12", "SAS", 125, 120, 245)
)

$r = 2
foreach ($game in $teamGames) {
    Set-TextValue $teamSheet.Cells.Item($r, 1) $game[0]
    for ($c = 1; $c -lt $game.Length; $c++) {
        $teamSheet.Cells.Item($r, $c + 1).Value = $game[$c]
    }
    $r += 1
}
